# Updates betexplorer turkey super-lig 2023-2024 sheet:
#  - several existing rows (which share the same kickoff date/time in column E)
#    had their match-detail columns (F:V) re-shuffled between rows
#  - four new match rows (174-177) are appended at the end

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Swap-Rows($r1, $r2) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

function Rotate-Rows($rowList) {
    # shifts values so that new[i] = old[i-1] (last wraps to first)
    $snapshots = @{}
    foreach ($r in $rowList) {
        $rowData = @{}
        foreach ($col in $cols) {
            $rowData[$col] = $ws.Range("$col$r").Value2
        }
        $snapshots[$r] = $rowData
    }
    $n = $rowList.Count
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $rowList[$i]
        $srcRow = $rowList[($i + $n - 1) % $n]
        $srcData = $snapshots[$srcRow]
        foreach ($col in $cols) {
            $ws.Range("$col$destRow").Value2 = $srcData[$col]
        }
    }
}

# Simple pairwise swaps (rows share identical kickoff date/time in column E)
Swap-Rows 89 90
Swap-Rows 91 92
Swap-Rows 119 120
Swap-Rows 130 131
Swap-Rows 132 133
Swap-Rows 155 156
Swap-Rows 167 168

# Three-way rotation: row97 <- old row99, row98 <- old row97, row99 <- old row98
Rotate-Rows @(97, 98, 99)

# Append four new match rows at the end of the table
$newRows = @(
    @{ Row = 174; A = 173; E = 45297.47916666666; F = "Rizespor";    G = 2; H = "Hatayspor";       I = 0;
       J = 2.24; K = "28/12/2024 19:12"; L = 1.88; M = "06/01/2024 11:22";
       N = 3.5;  O = "28/12/2024 19:12"; P = 3.57; Q = "06/01/2024 11:28";
       R = 3.25; S = "28/12/2024 19:12"; T = 4.62; U = "06/01/2024 11:22";
       V = "https://www.betexplorer.com/football/turkey/super-lig/rizespor-hatayspor/4UOCok15/" },
    @{ Row = 175; A = 174; E = 45297.58333333334; F = "Basaksehir";  G = 0; H = "Adana Demirspor"; I = 0;
       J = 2.62; K = "28/12/2024 19:12"; L = 1.98; M = "06/01/2024 13:58";
       N = 3.41; O = "28/12/2024 19:12"; P = 3.8;  Q = "06/01/2024 13:58";
       R = 2.75; S = "28/12/2024 19:12"; T = 3.83; U = "06/01/2024 13:59";
       V = "https://www.betexplorer.com/football/turkey/super-lig/basaksehir-adanademirspor/0xSekX1t/" },
    @{ Row = 176; A = 175; E = 45297.70833333334; F = "Ankaragucu";  G = 0; H = "Trabzonspor";      I = 1;
       J = 3.24; K = "28/12/2024 19:12"; L = 3.8;  M = "06/01/2024 16:59";
       N = 3.5;  O = "28/12/2024 19:12"; P = 3.34; Q = "06/01/2024 16:56";
       R = 2.25; S = "28/12/2024 19:12"; T = 2.15; U = "06/01/2024 16:56";
       V = "https://www.betexplorer.com/football/turkey/super-lig/ankaragucu-trabzonspor/jRK8nVoa/" },
    @{ Row = 177; A = 176; E = 45297.70833333334; F = "Samsunspor";  G = 1; H = "Karagumruk";       I = 0;
       J = 2.2;  K = "28/12/2024 19:12"; L = 2.47; M = "06/01/2024 16:58";
       N = 3.49; O = "28/12/2024 19:12"; P = 3.17; Q = "06/01/2024 16:58";
       R = 3.35; S = "28/12/2024 19:12"; T = 3.26; U = "06/01/2024 16:58";
       V = "https://www.betexplorer.com/football/turkey/super-lig/samsunspor-f-karagumruk/pQinVTGO/" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Range("A$r").Value2 = $nr.A
    $ws.Range("B$r").Value2 = "turkey"
    $ws.Range("C$r").Value2 = "super-lig"
    $ws.Range("D$r").Value2 = "2023-2024"
    $ws.Range("E$r").Value2 = $nr.E
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $nr[$col]
    }

    # Match the existing styling used by the rest of the table: column A uses
    # the bold/bordered/centered style, column E uses the datetime number format.
    $ws.Range("A173").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("E173").Copy() | Out-Null
    $ws.Range("E$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0
